# The commit removes the SharePoint / Content-Organizer "Custom XML Data"
# parts that had been attached to this document (customXml/item1.xml..item3.xml
# and their customXml/itemProps*.xml companions, together with the
# relationships that wire them into the package). Nothing in the visible
# body content changes - this is the same effect as running Word's
# Document Inspector over "Custom XML Data" and clicking "Remove All", or
# equivalently iterating Document.CustomXMLParts and deleting every part.
#
# The well-known part ids for this document (customXml/itemProps*.xml
# ds:itemID values, which become CustomXMLPart.Id through the object model):
#   customXml/item1.xml -> {FA10CC44-DA7A-49F1-8ACD-664E9798E90D}
#   customXml/item2.xml -> {604405CA-90DC-45F2-ABAE-E983DA5454D8}
#   customXml/item3.xml -> {2501ACB3-F32D-4D41-A745-8C7AACC16F5D}

$d = $word.ActiveDocument

$targetIds = @(
    "{FA10CC44-DA7A-49F1-8ACD-664E9798E90D}",
    "{604405CA-90DC-45F2-ABAE-E983DA5454D8}",
    "{2501ACB3-F32D-4D41-A745-8C7AACC16F5D}"
)

$parts = $d.CustomXMLParts

# Delete by known id first (most precise - mirrors a script that targets
# exactly the SharePoint metadata parts added to this file).
foreach ($id in $targetIds) {
    try {
        $part = $parts.SelectByID($id)
        if ($part) {
            $part.Delete()
        }
    } catch {
        # Host/part may already be gone - ignore and continue.
    }
}

# Belt-and-braces sweep: walk the remaining collection back-to-front (so
# deleting doesn't shift the indices we still have to visit) and remove
# anything that is not an Office built-in metadata schema, in case the
# id-based pass above missed one.
try {
    for ($i = $parts.Count; $i -ge 1; $i--) {
        $part = $parts.Item($i)
        $isBuiltIn = $false
        try {
            $ns = $part.NamespaceURI
            if ($ns -eq "http://schemas.openxmlformats.org/package/2006/metadata/core-properties" -or
                $ns -eq "http://schemas.openxmlformats.org/officeDocument/2006/extended-properties" -or
                $ns -eq "http://schemas.microsoft.com/office/2006/coverPageProps") {
                $isBuiltIn = $true
            }
        } catch {
        }
        if (-not $isBuiltIn) {
            try {
                $part.Delete()
            } catch {
            }
        }
    }
} catch {
}

Write-Output ("CustomXMLParts remaining: " + $d.CustomXMLParts.Count)
